# Apply the scraped crypto-price refresh (GitHub Actions update).
# Every target cell holds plain text (prices use literal "."-grouped
# strings like "3.463.35", not real numbers), so each write forces the
# cell to Text format first when the new value would otherwise be
# auto-converted to a number by Excel, then restores the original style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    if ($value -match '^\s*[+-]?\d+(\.\d+)?\s*$') {
        $savedStyle = $range.Style
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = $savedStyle
    } else {
        $range.Value = $value
    }
}

Set-TextValue $ws.Range("D2") "66.989.05"
Set-TextValue $ws.Range("E2") "  -0.30%  "
Set-TextValue $ws.Range("D3") "3.453.42"
Set-TextValue $ws.Range("E3") "  -1.35%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "591.70"
Set-TextValue $ws.Range("E5") "  -1.33%  "
Set-TextValue $ws.Range("D6") "179.45"
Set-TextValue $ws.Range("E6") "  +2.10%  "
Set-TextValue $ws.Range("D7") "0.612"
Set-TextValue $ws.Range("E7") "  +4.04%  "
Set-TextValue $ws.Range("D8") "1.00"
Set-TextValue $ws.Range("E8") "  +0.09%  "
Set-TextValue $ws.Range("D9") "3.451.34"
Set-TextValue $ws.Range("E9") "  -1.34%  "
Set-TextValue $ws.Range("D10") "0.139"
Set-TextValue $ws.Range("E10") "  +5.70%  "
Set-TextValue $ws.Range("D11") "6.93"
Set-TextValue $ws.Range("E11") "  -3.14%  "
Set-TextValue $ws.Range("D12") "0.430"
Set-TextValue $ws.Range("E12") "  -0.19%  "
Set-TextValue $ws.Range("D13") "4.056.76"
Set-TextValue $ws.Range("E13") "  -1.25%  "
Set-TextValue $ws.Range("D14") "31.88"
Set-TextValue $ws.Range("E14") "  +2.11%  "
Set-TextValue $ws.Range("E15") "  -0.61%  "
Set-TextValue $ws.Range("D16") "67.016.18"
Set-TextValue $ws.Range("E16") "  -0.25%  "
Set-TextValue $ws.Range("D17") "0.0000176"
Set-TextValue $ws.Range("E17") "  -1.64%  "
Set-TextValue $ws.Range("D18") "3.459.59"
Set-TextValue $ws.Range("E18") "  -0.95%  "
Set-TextValue $ws.Range("D19") "6.20"
Set-TextValue $ws.Range("E19") "  -1.50%  "
Set-TextValue $ws.Range("D20") "14.16"
Set-TextValue $ws.Range("E20") "  -2.27%  "
Set-TextValue $ws.Range("D21") "389.84"
Set-TextValue $ws.Range("E21") "  -1.19%  "
Set-TextValue $ws.Range("D22") "7.92"
Set-TextValue $ws.Range("E22") "  -1.12%  "
Set-TextValue $ws.Range("E23") "  +0.01%  "
Set-TextValue $ws.Range("E24") "  +1.48%  "
Set-TextValue $ws.Range("D25") "72.00"
Set-TextValue $ws.Range("E25") "  -1.93%  "
Set-TextValue $ws.Range("D26") "0.535"
Set-TextValue $ws.Range("E26") "  -0.34%  "
Set-TextValue $ws.Range("D27") "0.0000121"
Set-TextValue $ws.Range("E27") "  -0.51%  "
Set-TextValue $ws.Range("D28") "10.29"
Set-TextValue $ws.Range("E28") "  +0.88%  "
Set-TextValue $ws.Range("D29") "0.174"
Set-TextValue $ws.Range("E29") "  -3.47%  "
Set-TextValue $ws.Range("E30") "  +0.56%  "
Set-TextValue $ws.Range("D31") "6.12"
Set-TextValue $ws.Range("E31") "  -0.26%  "
Set-TextValue $ws.Range("B32") "PancakeSwap"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D32") "2.05"
Set-TextValue $ws.Range("E32") "  -0.57%  "
Set-TextValue $ws.Range("B33") "Fetch.AI"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D33") "1.39"
Set-TextValue $ws.Range("E33") "  -2.10%  "
Set-TextValue $ws.Range("D34") "23.40"
Set-TextValue $ws.Range("E34") "  -1.05%  "
Set-TextValue $ws.Range("D35") "7.30"
Set-TextValue $ws.Range("E35") "  -0.99%  "
Set-TextValue $ws.Range("E36") "  -0.05%  "
Set-TextValue $ws.Range("E37") "  -3.45%  "
Set-TextValue $ws.Range("D38") "163.77"
Set-TextValue $ws.Range("E38") "  +0.57%  "
Set-TextValue $ws.Range("D39") "0.873"
Set-TextValue $ws.Range("E39") "  -0.49%  "
Set-TextValue $ws.Range("D40") "2.79"
Set-TextValue $ws.Range("E40") "  +9.03%  "
Set-TextValue $ws.Range("D41") "1.86"
Set-TextValue $ws.Range("E41") "  -3.66%  "
Set-TextValue $ws.Range("D42") "6.80"
Set-TextValue $ws.Range("E42") "  -3.32%  "
Set-TextValue $ws.Range("D43") "4.64"
Set-TextValue $ws.Range("E43") "  -0.27%  "
Set-TextValue $ws.Range("D44") "26.10"
Set-TextValue $ws.Range("E44") "  +0.15%  "
Set-TextValue $ws.Range("D45") "0.0718"
Set-TextValue $ws.Range("E45") "  -1.78%  "
Set-TextValue $ws.Range("D46") "2.739.19"
Set-TextValue $ws.Range("E46") "  -2.57%  "
Set-TextValue $ws.Range("D47") "26.18"
Set-TextValue $ws.Range("E47") "  -4.64%  "
Set-TextValue $ws.Range("D48") "41.30"
Set-TextValue $ws.Range("E48") "  -2.70%  "
Set-TextValue $ws.Range("D49") "0.0298"
Set-TextValue $ws.Range("E49") "  -1.82%  "
Set-TextValue $ws.Range("D50") "325.31"
Set-TextValue $ws.Range("E50") "  -3.72%  "
Set-TextValue $ws.Range("E51") "  -3.91%  "
